# #318 [done]: [neon-cli]: Rename [dns-hosts] to [hive dns]
#
# Fix the typo "certificated" -> "certificate" in the paragraph about
# operators using a real trusted certificate authority. Word records the
# cursor position of the last text change in a hidden "_GoBack" bookmark;
# since a document can only have one bookmark of a given name, re-adding
# "_GoBack" at the edit site automatically relocates it away from its
# previous position (around "sysadmin/password").

$d = $word.ActiveDocument

# Locate the sentence containing the typo and capture its Range.
$editRange = $d.Content
$found = $editRange.Find.Execute(
    "certificated signed by a real trusted 3",  # FindText
    $true,                                      # MatchCase
    $false,                                     # MatchWholeWord
    $false,                                     # MatchWildcards
    $false,                                     # MatchSoundsLike
    $false,                                     # MatchAllWordForms
    $true,                                      # Forward
    1,                                          # Wrap (wdFindContinue)
    $false,                                     # Format
    "",                                         # ReplaceWith
    0                                           # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find the 'certificated' sentence to fix."
}

# "certificate" is 11 characters; the very next character is the stray "d"
# that needs to be removed ("certificated" -> "certificate").
$certificateEnd = $editRange.Start + 11
$strayD = $d.Range($certificateEnd, $certificateEnd + 1)
if ($strayD.Text -ne "d") {
    throw "Unexpected character where 'd' was expected: [$($strayD.Text)]"
}
$strayD.Text = ""

# Drop a (collapsed) "_GoBack" bookmark right where the edit happened, as
# Word itself would after typing here. Re-using the name moves it off the
# old "sysadmin/password" location automatically.
$bookmarkSpot = $d.Range($certificateEnd, $certificateEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

Write-Output "Fixed 'certificated' -> 'certificate' and relocated the _GoBack bookmark."
